$d = $word.ActiveDocument

# Locate the paragraph that contains "45" (the end of "12345") and
# position right after it, before the trailing bookmark.
$target = $d.Paragraphs(2).Range
$find = $target.Find
$find.Execute("45", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Collapse to the end of "45" so we insert right after it.
$insertionPoint = $find.Parent
$insertionPoint.Collapse(0)

# Insert a paragraph break followed by the new text "Good Bay".
$insertionPoint.InsertParagraphAfter()
$insertionPoint.Collapse(0)
$insertionPoint.Move(4, 1) | Out-Null
$insertionPoint.Text = "Good Bay"
